# Daily attendance processing - 2025-10-30 19:42:16
#
# Column G ("Recorded By") holds a comma-separated list of the
# users/processes that touched each attendance row. Re-processing the
# daily log rotates that list so the first entry moves to the end
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Cells with only a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -notlike "*,*") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -le 1) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value2 = $rotated
}
